$p = $ppt.ActivePresentation

$newDate = "2025/7/2"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }

            if ($isDatePlaceholder) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -ne $newDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Update the slide master's date & time placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Update every slide layout's date & time placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
